# Apply the "Filled in example manifests" commit:
#   1. README!C2 - reword the checklist description.
#   2. HiddenDropdowns!F80 - add a new tissue/body-site value ("Inflorescence")
#      to the end of the hidden dropdown source list, and extend the
#      "sample" sheet's data-validation range that points at it.

$wb = $excel.ActiveWorkbook

# --- 1. README description text -------------------------------------------
$readme = $wb.Worksheets.Item("README")
$readme.Unprotect()
$readme.Range("C2").Value = "Spatial transcriptomics sequencing, using  Minimum Information about any (x) Sequence (MIxS) standard, detailing the contextual information about sampling and sequencing processes."
$readme.Columns.Item(3).AutoFit()
$readme.Protect()

# --- 2. Add "Inflorescence" to the hidden tissue dropdown list -------------
$hidden = $wb.Worksheets.Item("HiddenDropdowns")
$hidden.Range("F80").Value = "Inflorescence"

# --- 3. Point the "sample" sheet's tissue dropdown at the new, larger range
$sample = $wb.Worksheets.Item("sample")
$sample.Unprotect()
$sample.Range("F5:F1005").Validation.Formula1 = "HiddenDropdowns!`$F`$5:`$F`$80"
$sample.Protect()
